$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated county statistics figures ---
$ws.Range("D12").Value = 480
$ws.Range("E12").Value = 94
$ws.Range("H12").Value = 0.83623693379790942
$ws.Range("I12").Value = 0.16376306620209058
$ws.Range("K12").Value = 40912
$ws.Range("L12").Value = 2949
$ws.Range("O12").Value = 0.93276487084197812
$ws.Range("P12").Value = 0.06723512915802193346
$ws.Range("E14").Value = 188
$ws.Range("F14").Value = 1106
$ws.Range("H14").Value = 0.83001808318264014
$ws.Range("I14").Value = 0.16998191681735986
$ws.Range("L14").Value = 7531
$ws.Range("M14").Value = 92294
$ws.Range("O14").Value = 0.91840206297267424
$ws.Range("P14").Value = 0.0815979370273257204
$ws.Range("D21").Value = 1039
$ws.Range("E21").Value = 251
$ws.Range("H21").Value = 0.8054263565891473
$ws.Range("I21").Value = 0.1945736434108527
$ws.Range("K21").Value = 92929
$ws.Range("L21").Value = 8554
$ws.Range("O21").Value = 0.91571002039750504
$ws.Range("P21").Value = 0.08428997960249499932
$ws.Range("D22").Value = 1017
$ws.Range("E22").Value = 162
$ws.Range("H22").Value = 0.86259541984732824
$ws.Range("I22").Value = 0.13740458015267176
$ws.Range("K22").Value = 89974
$ws.Range("L22").Value = 5638
$ws.Range("O22").Value = 0.94103250637995228
$ws.Range("P22").Value = 0.05896749362004768957
$ws.Range("D23").Value = 1186
$ws.Range("E23").Value = 534
$ws.Range("H23").Value = 0.68953488372093019
$ws.Range("I23").Value = 0.31046511627906975
$ws.Range("K23").Value = 105564
$ws.Range("L23").Value = 25495
$ws.Range("O23").Value = 0.80546929245606935
$ws.Range("P23").Value = 0.19453070754393059
$ws.Range("D28").Value = 1762
$ws.Range("E28").Value = 462
$ws.Range("H28").Value = 0.79226618705035967
$ws.Range("I28").Value = 0.2077338129496403
$ws.Range("K28").Value = 155146
$ws.Range("L28").Value = 16188
$ws.Range("O28").Value = 0.9055178773623449
$ws.Range("P28").Value = 0.09448212263765511099
$ws.Range("D29").Value = 2670
$ws.Range("E29").Value = 442
$ws.Range("F29").Value = 3112
$ws.Range("H29").Value = 0.85796915167095111
$ws.Range("I29").Value = 0.14203084832904883
$ws.Range("K29").Value = 248104
$ws.Range("L29").Value = 20402
$ws.Range("M29").Value = 268506
$ws.Range("O29").Value = 0.92401659553231585
$ws.Range("P29").Value = 0.0759834044676841458
$ws.Range("D32").Value = 365
$ws.Range("E32").Value = 162
$ws.Range("H32").Value = 0.69259962049335866
$ws.Range("I32").Value = 0.30740037950664134
$ws.Range("K32").Value = 36846
$ws.Range("L32").Value = 7981
$ws.Range("O32").Value = 0.82195997947665467
$ws.Range("P32").Value = 0.1780400205233453
$ws.Range("D34").Value = 926
$ws.Range("E34").Value = 215
$ws.Range("H34").Value = 0.81156879929886061
$ws.Range("I34").Value = 0.18843120070113936
$ws.Range("K34").Value = 79708
$ws.Range("L34").Value = 7730
$ws.Range("O34").Value = 0.91159450124659758
$ws.Range("P34").Value = 0.08840549875340240893
$ws.Range("D35").Value = 656
$ws.Range("E35").Value = 128
$ws.Range("H35").Value = 0.83673469387755106
$ws.Range("I35").Value = 0.16326530612244897
$ws.Range("K35").Value = 58016
$ws.Range("L35").Value = 4575
$ws.Range("O35").Value = 0.92690642424629743
$ws.Range("P35").Value = 0.07309357575370260818
$ws.Range("D38").Value = 582
$ws.Range("E38").Value = 178
$ws.Range("H38").Value = 0.76578947368421058
$ws.Range("I38").Value = 0.23421052631578948
$ws.Range("K38").Value = 50599
$ws.Range("L38").Value = 4685
$ws.Range("O38").Value = 0.9152557702047609
$ws.Range("P38").Value = 0.0847442297952391288
$ws.Range("D43").Value = 719
$ws.Range("E43").Value = 118
$ws.Range("H43").Value = 0.85902031063321382
$ws.Range("I43").Value = 0.14097968936678615
$ws.Range("K43").Value = 60933
$ws.Range("L43").Value = 3061
$ws.Range("O43").Value = 0.95216739069287748
$ws.Range("P43").Value = 0.04783260930712254083
$ws.Range("D44").Value = 1438
$ws.Range("E44").Value = 1078
$ws.Range("H44").Value = 0.57154213036565982
$ws.Range("I44").Value = 0.42845786963434024
$ws.Range("K44").Value = 131185
$ws.Range("L44").Value = 54445
$ws.Range("O44").Value = 0.70670150298981849
$ws.Range("P44").Value = 0.29329849701018157
$ws.Range("D45").Value = 1655
$ws.Range("E45").Value = 444
$ws.Range("H45").Value = 0.78847070033349209
$ws.Range("I45").Value = 0.21152929966650785
$ws.Range("K45").Value = 144014
$ws.Range("L45").Value = 17360
$ws.Range("O45").Value = 0.89242381052709852
$ws.Range("P45").Value = 0.10757618947290146
$ws.Range("D48").Value = 724
$ws.Range("E48").Value = 122
$ws.Range("H48").Value = 0.85579196217494091
$ws.Range("I48").Value = 0.14420803782505912
$ws.Range("K48").Value = 63631
$ws.Range("L48").Value = 3733
$ws.Range("O48").Value = 0.94458464461730296
$ws.Range("P48").Value = 0.05541535538269699163
$ws.Range("D54").Value = 36853
$ws.Range("E54").Value = 10260
$ws.Range("H54").Value = 0.78222571264831364
$ws.Range("I54").Value = 0.21777428735168636
$ws.Range("K54").Value = 3308839
$ws.Range("L54").Value = 431875
$ws.Range("O54").Value = 0.8845474420124072
$ws.Range("P54").Value = 0.11545255798759274
$ws.Range("D87").Value = 68
$ws.Range("E87").Value = 45
$ws.Range("H87").Value = 0.60176991150442483
$ws.Range("I87").Value = 0.39823008849557523
$ws.Range("K87").Value = 5944
$ws.Range("L87").Value = 2385
$ws.Range("O87").Value = 0.71365109857125708
$ws.Range("P87").Value = 0.28634890142874292
$ws.Range("D88").Value = 1650
$ws.Range("E88").Value = 499
$ws.Range("H88").Value = 0.76779897626803162
$ws.Range("I88").Value = 0.23220102373196835
$ws.Range("K88").Value = 137413
$ws.Range("L88").Value = 22184
$ws.Range("O88").Value = 0.86099989348170702
$ws.Range("P88").Value = 0.13900010651829295
$ws.Range("D139").Value = 52156
$ws.Range("E139").Value = 13857
$ws.Range("H139").Value = 0.7900868010846348
$ws.Range("I139").Value = 0.20991319891536517
$ws.Range("K139").Value = 4689883
$ws.Range("L139").Value = 586054
$ws.Range("O139").Value = 0.88891944691530622
$ws.Range("P139").Value = 0.11108055308469376

# --- Column width tweaks ---
$ws.Columns.Item(4).ColumnWidth = 8.7109375   # column D
$ws.Columns.Item(11).ColumnWidth = 8.5703125  # column K

# --- Selection state ---
$ws.Range("C4").Select()
